$wb = $excel.ActiveWorkbook

# --- Language / content fixes on the Banner_Text sheet ---
$banner = $wb.Worksheets.Item("Banner_Text")
$banner.Range("B2").Value = "Sri Brahmatantra Svatantra Parakāla Swāmi Mutt Guru Paramparā"
$banner.Range("B4").Value = "Sri Brahmatantra Svatantra Parakāla Swāmi Mutt Ācāryas"
$banner.Range("B6").Value = "Sri Parakāla Swāmi Mutt – The Eternal Lineage of Sri Vedānta Deśika"

# --- Make Banner_Text the active/selected sheet (was acharyan_captions) ---
$banner.Activate() | Out-Null
$banner.Range("B21").Select() | Out-Null
